$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KB")

# Insert a new column before column B so that layout becomes:
# A=min_pH, B=max_pH, C=N, D=P, E=K
$ws.Range("B1").EntireColumn.Insert()

# Header row (set max_pH before min_pH so shared-string order matches target)
$ws.Range("B1").Value = "max_pH"
$ws.Range("A1").Value = "min_pH"
$ws.Range("C1").Value = "N"
$ws.Range("D1").Value = "P"
$ws.Range("E1").Value = "K"

# Data rows: A (min_pH) already has old pH values (0, 0.5, 1, ... 14)
# Fill B (max_pH) = A + 0.49 (except row 6 has floating point quirk 2.4900000000000002)
$maxVals = @(0.49,0.99,1.49,1.99,2.4900000000000002,2.99,3.49,3.99,4.49,4.99,5.49,5.99,6.49,6.99,7.49,7.99,8.49,8.99,9.49,9.99,10.49,10.99,11.49,11.99,12.49,12.99,13.49,13.99,14.49)

for ($i = 0; $i -lt $maxVals.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $maxVals[$i]
}

# N, P, K values for rows 2..30 (low/mid/high)
$nVals = @("low","low","low","low","low","low","low","low","low","low","mid","mid","high","high","high","high","high","mid","mid","low","low","low","low","low","low","low","low","low","low")
$pVals = @("low","low","low","low","low","low","low","low","low","low","low","low","mid","high","high","high","mid","low","high","high","high","high","high","high","high","high","high","high","high")
$kVals = @("low","low","low","low","low","low","low","low","low","low","mid","mid","high","high","high","high","high","high","high","high","high","high","high","high","high","high","high","high","high")

for ($i = 0; $i -lt $nVals.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $nVals[$i]
    $ws.Cells.Item($row, 4).Value = $pVals[$i]
    $ws.Cells.Item($row, 5).Value = $kVals[$i]
}

$ws.Range("J22").Select()
